# This workbook is a small single-sheet contact list:
#   A1 = phone number
#   A2 = "Sumit"
#   A3 = email-ish text (becomes a mailto hyperlink)
#   A4 = "sum@gmail.com" (existing mailto hyperlink)
#   A5 = "Sumit Gokhe"
#   A6 = phone number
#
# The edit:
#   - A1 / A6 (the two phone numbers) get an explicit left-aligned style
#   - A3's text changes from "sum@7889" to "sum@gmail.c" and becomes a
#     mailto hyperlink (styled like the existing A4 hyperlink)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Left-align the two numeric/phone cells (creates the new cellXfs entry
# with horizontal="left").
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A6").HorizontalAlignment = -4131

# Update A3's text and turn it into a mailto hyperlink, matching the
# existing hyperlink on A4.
$ws.Range("A3").Value = "sum@gmail.c"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sum@gmail.c")
$ws.Range("A3").Style = "Hyperlink"
